# Update DataBI\dollar_to_yuan.xlsx with a new conversion data row (row 5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new cells to be stored as text (matching the existing inline/shared
# string cells above) rather than being auto-detected as a date / number.
$ws.Range("A5:B5").NumberFormat = "@"
$ws.Range("A5").Value = "2024-12-16"
$ws.Range("B5").Value = "7.283981"

# Drop the temporary text format again so the new cells end up with the same
# (default) style as the rest of the sheet, just like the original rows.
$ws.Range("A5:B5").ClearFormats()
